$wb = $excel.ActiveWorkbook

# --- Rename sheets (new task-order run ids) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16512556291226792"
$wb.Worksheets.Item(2).Name = "NB_TO-16512556336367197"
$wb.Worksheets.Item(3).Name = "RS_TO-16512556336437001"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512556337015522"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512556337776427"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1651255629085417.csv"
$ws1.Range("B3").Value = "GNG_stims-1651255629106341.csv"
$ws1.Range("B4").Value = "go_stims-16512556291080334.csv"
$ws1.Range("B5").Value = "GNG_stims-16512556291216798.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1651255632866564.csv"
$ws2.Range("B3").Value = "OB-16512556310413184.csv"
$ws2.Range("B4").Value = "TB-16512556330822542.csv"
$ws2.Range("B5").Value = "ZB-match_0-16512556292675722.csv"
$ws2.Range("B6").Value = "OB-1651255630673402.csv"
$ws2.Range("B7").Value = "TB-1651255633621946.csv"
$ws2.Range("B8").Value = "OB-165125563003412.csv"
$ws2.Range("B9").Value = "ZB-match_1-1651255629568353.csv"
$ws2.Range("B10").Value = "ZB-match_7-16512556293940644.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512556336681576.csv"
$ws4.Range("B3").Value = "ZM_stims-16512556336456974.csv"
$ws4.Range("B4").Value = "MM_stims-16512556336836991.csv"
$ws4.Range("B5").Value = "ZM_stims-16512556336696982.csv"
$ws4.Range("B6").Value = "MM_stims-16512556336995518.csv"
$ws4.Range("B7").Value = "ZM_stims-16512556336847005.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1651255633746462.csv"
$ws5.Range("B3").Value = "SAT_stims-16512556337312913.csv"
$ws5.Range("B4").Value = "SAT_stims-16512556337065053.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651255633762201.csv"
